$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150-216 down to 151-217
$ws.Rows("150:150").Insert()

# Populate the new row 150 with the values from the commit (new weekly entry)
$ws.Range("A150").Value = 11
$ws.Range("B150").Value = "Vega Monumental Concepción"
$ws.Range("C150").Value = "Bíobío"
$ws.Range("D150").Value = 45119
$ws.Range("D150").NumberFormat = $ws.Range("D151").NumberFormat
$ws.Range("E150").Value = 8
$ws.Range("F150").Value = 100112021
$ws.Range("G150").Value = "Ají"
$ws.Range("H150").Value = "Americana (o)"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 50
$ws.Range("K150").Value = 28000
$ws.Range("L150").Value = 30000
$ws.Range("M150").Value = 29200
$ws.Range("N150").Value = "$/caja 25 kilos"
$ws.Range("O150").Value = "Provincia de Limarí"
$ws.Range("P150").Value = 1168
$ws.Range("Q150").Value = 25
$ws.Range("R150").Value = "Hortaliza"
